$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Import" cell to drop FlaggedPatient in favor of PatientFlag
$ws.Range("C2").Value = "org.openmrs.Patient,org.openmrs.module.patientflags.PatientFlag,org.openmrs.module.drools.calculation.Operator, static org.openmrs.module.drools.utils.DroolsDateUtils.daysAgo"

# Update the rule table's CONDITION / ACTION cells for "No Existing Sepsis Flag" / "Has Sepsis Flag With Message" / "Flag Patient"
$ws.Range("H11").Value = "not PatientFlag(patient == `$patient)"
$ws.Range("I11").Value = "`$flag: PatientFlag(patient == `$patient, message == `$param)"
$ws.Range("J11").Value = "insert(new PatientFlag(`$patient, null, `$param));"

# Update the view state: scroll position and active selected cell
$ws.Range("J11").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.DisplayGridlines = $true
